# Summary.xlsx RAD data refresh
# - Updates the "Date" column (B) timestamps on several sheets to reflect a
#   newer Katalon test run (Wed Oct 11 2023 afternoon), replacing the old
#   Fri Oct 06 / Mon Oct 09 runs.
# - Restores the active sheet's (Existing) selection to B17.
#
# Cell values are assigned in the same relative order the new run's shared
# strings appear in the source workbook (sheet-by-sheet, top-to-bottom) so the
# resulting shared-string table is built up consistently.

$wb = $excel.ActiveWorkbook

# --- Estimated ---
$ws = $wb.Worksheets.Item("Estimated")
$ws.Range("B2").Value = "Wed Oct 11 13:18:24 EDT 2023"
$ws.Range("B3").Value = "Wed Oct 11 13:19:04 EDT 2023"
$ws.Range("B4").Value = "Wed Oct 11 13:19:42 EDT 2023"
$ws.Range("B5").Value = "Wed Oct 11 13:20:19 EDT 2023"
$ws.Range("B6").Value = "Wed Oct 11 13:20:56 EDT 2023"
$ws.Range("B7").Value = "Wed Oct 11 13:21:33 EDT 2023"

# --- Existing ---
$ws = $wb.Worksheets.Item("Existing")
$ws.Range("B2").Value = "Wed Oct 11 13:02:26 EDT 2023"
$ws.Range("B3").Value = "Wed Oct 11 13:03:07 EDT 2023"
$ws.Range("B4").Value = "Wed Oct 11 13:03:45 EDT 2023"
$ws.Range("B5").Value = "Wed Oct 11 13:04:23 EDT 2023"
$ws.Range("B6").Value = "Wed Oct 11 13:05:01 EDT 2023"
$ws.Range("B7").Value = "Wed Oct 11 13:05:38 EDT 2023"
$ws.Range("B8").Value = "Wed Oct 11 13:06:17 EDT 2023"
$ws.Range("B9").Value = "Wed Oct 11 13:06:56 EDT 2023"
$ws.Range("B10").Value = "Wed Oct 11 13:07:36 EDT 2023"
$ws.Range("B11").Value = "Wed Oct 11 13:08:15 EDT 2023"
$ws.Range("B12").Value = "Wed Oct 11 13:08:52 EDT 2023"

# --- Extension ---
$ws = $wb.Worksheets.Item("Extension")
$ws.Range("B2").Value = "Wed Oct 11 13:25:48 EDT 2023"
$ws.Range("B3").Value = "Wed Oct 11 13:26:28 EDT 2023"
$ws.Range("B4").Value = "Wed Oct 11 13:27:04 EDT 2023"
$ws.Range("B5").Value = "Wed Oct 11 13:27:39 EDT 2023"
$ws.Range("B6").Value = "Wed Oct 11 13:28:15 EDT 2023"
$ws.Range("B7").Value = "Wed Oct 11 13:28:51 EDT 2023"

# --- NewTaxReturn ---
$ws = $wb.Worksheets.Item("NewTaxReturn")
$ws.Range("B2").Value = "Wed Oct 11 16:17:53 EDT 2023"
$ws.Range("B3").Value = "Wed Oct 11 16:18:32 EDT 2023"
$ws.Range("B4").Value = "Wed Oct 11 16:19:08 EDT 2023"
$ws.Range("B5").Value = "Wed Oct 11 16:19:43 EDT 2023"
$ws.Range("B6").Value = "Wed Oct 11 16:20:19 EDT 2023"
$ws.Range("B7").Value = "Wed Oct 11 16:20:55 EDT 2023"
$ws.Range("B8").Value = "Wed Oct 11 16:21:31 EDT 2023"
$ws.Range("B9").Value = "Wed Oct 11 16:22:06 EDT 2023"
$ws.Range("B10").Value = "Wed Oct 11 16:22:42 EDT 2023"
$ws.Range("B11").Value = "Wed Oct 11 16:23:18 EDT 2023"
$ws.Range("B12").Value = "Wed Oct 11 16:23:54 EDT 2023"
$ws.Range("B13").Value = "Wed Oct 11 16:24:30 EDT 2023"
$ws.Range("B14").Value = "Wed Oct 11 16:25:05 EDT 2023"
$ws.Range("B15").Value = "Wed Oct 11 16:25:40 EDT 2023"
$ws.Range("B16").Value = "Wed Oct 11 16:26:16 EDT 2023"

# --- Personal_EL ---
$ws = $wb.Worksheets.Item("Personal_EL")
$ws.Range("B2").Value = "Wed Oct 11 16:34:18 EDT 2023"

# --- Personal_IND ---
$ws = $wb.Worksheets.Item("Personal_IND")
$ws.Range("B2").Value = "Wed Oct 11 16:42:56 EDT 2023"
$ws.Range("B3").Value = "Wed Oct 11 16:43:32 EDT 2023"
$ws.Range("B4").Value = "Wed Oct 11 16:44:06 EDT 2023"
$ws.Range("B5").Value = "Wed Oct 11 16:44:40 EDT 2023"
$ws.Range("B6").Value = "Wed Oct 11 16:45:14 EDT 2023"

# --- Personal_JNT ---
$ws = $wb.Worksheets.Item("Personal_JNT")
$ws.Range("B2").Value = "Wed Oct 11 17:06:12 EDT 2023"
$ws.Range("B3").Value = "Wed Oct 11 17:06:56 EDT 2023"
$ws.Range("B4").Value = "Wed Oct 11 17:07:38 EDT 2023"
$ws.Range("B5").Value = "Wed Oct 11 17:08:19 EDT 2023"
$ws.Range("B6").Value = "Wed Oct 11 17:09:01 EDT 2023"

# --- Restore the active sheet ("Existing") and its selection (B17) ---
$ws = $wb.Worksheets.Item("Existing")
$ws.Activate()
$ws.Range("B17").Select()
